# Update the "想去人数" (F column) counts that changed between crawls.
# The same set of rows/values needs updating on both the "展览" sheet
# and the "全部类型" sheet (they mirror the same underlying data).

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 10185
    8  = 116
    12 = 434
    15 = 522
    16 = 17
    19 = 321
    23 = 35
    25 = 73
    26 = 770
    27 = 1347
    31 = 53
    33 = 29
    36 = 183
    37 = 220
    41 = 98
    43 = 546
    44 = 3148
    46 = 156
    47 = 819
    49 = 35
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
